# Time Record Manager - 21 Jan 2025
# Update the two Global-Search-User project entries on the TimeClockRecorder
# sheet, and leave the workbook focused on that sheet (matching the author's
# last-saved selection) instead of Error_Message.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TimeClockRecorder")

$ws.Range("A2").Value = "Project Beach - Tax-StoicLane-FVA-122317"
$ws.Range("A3").Value = "Wittliff_Goodkind-Wittliff Cutter PLLC-FVA-109581"

# Make TimeClockRecorder the active sheet/tab, with A10 selected - this moves
# tabSelected off Error_Message (previously active) and onto this sheet.
$ws.Activate()
$ws.Range("A10").Select()
